$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column H header "Who" (bold like the other headers A1:G1)
$ws.Range("H1").Value = "Who"
$ws.Range("H1").Font.Bold = $true

# Initials for who ran each experiment
$ws.Range("H2").Value = "EZ"
$ws.Range("H3").Value = "EZ"
$ws.Range("H4").Value = "CF"
$ws.Range("H5").Value = "CF"
$ws.Range("H6").Value = "CF"
$ws.Range("H7").Value = "CF"
$ws.Range("H8").Value = "CF"

# Update the active selection to I10
$ws.Range("I10").Select() | Out-Null
